$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that get shuffled/updated by this weekly data refresh.
$cols = @("D", "H", "I", "J", "K", "L", "M", "P")

# Capture the "before" values for every relevant column across rows 2-18
# so we can re-distribute them according to the new row order.
$before = @{}
foreach ($c in $cols) {
    $before[$c] = @{}
    for ($r = 2; $r -le 18; $r++) {
        $before[$c][$r] = $ws.Range("$c$r").Value2
    }
}

# Mapping: new row -> old row whose data it now carries.
$rowMap = @{
    2  = 5
    3  = 13
    4  = 9
    5  = 14
    6  = 16
    7  = 8
    8  = 6
    9  = 11
    10 = 2
    11 = 15
    12 = 3
    13 = 12
    14 = 18
    15 = 10
    16 = 17
    17 = 7
    18 = 4
}

foreach ($r in 2..18) {
    $src = $rowMap[$r]
    foreach ($c in $cols) {
        $ws.Range("$c$r").Value = $before[$c][$src]
    }
}
